# LMR_Table.xlsx — "Implemented SEE - testing needed"
#
# The original sole sheet ("Tabelle1", the big Late-Move-Reduction grid)
# is renamed to "LMR". A brand new sheet, itself named "Tabelle1", is
# inserted right after it and holds a small SEE (Static Exchange
# Evaluation) "Pow"/"Base" helper calculation: column D holds an
# incrementing Offset (1..32), column E computes
#   = $B$1 + ((D<n> + $B$2) ^ $B$3)
# with B1=Offset(4), B2=Base(1), B3=Pow(1.5), formatted as an integer.

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet, add the new one right after it -----------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LMR"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle1"

# --- header / parameter cells ---------------------------------------------
# Insertion order matters for the shared-string table layout: "Pow" is
# typed before "Base".
$ws2.Range("A1").Value = "Offset"
$ws2.Cells.Item(1, 2).Value = 4
$ws2.Cells.Item(1, 4).Value = 1

$ws2.Range("A3").Value = "Pow"
$ws2.Cells.Item(3, 2).Value = 1.5
$ws2.Cells.Item(3, 4).Value = 3

$ws2.Range("A2").Value = "Base"
$ws2.Cells.Item(2, 2).Value = 1
$ws2.Cells.Item(2, 4).Value = 2

# --- offsets 4..32 in column D --------------------------------------------
for ($r = 4; $r -le 32; $r++) {
    $ws2.Cells.Item($r, 4).Value = $r
}

# --- column E: the SEE formula, shared down the whole range --------------
$ws2.Range("E1").Formula = "=`$B`$1+((D1+`$B`$2)^`$B`$3)"
$ws2.Range("E2:E32").Formula = "=`$B`$1+((D2+`$B`$2)^`$B`$3)"
$ws2.Range("E1:E32").NumberFormat = "0"

$ws2.Columns.Item(5).ColumnWidth = 10.25

# match the workbook's usual 2cm top/bottom page margins (inches, points API)
$ws2.PageSetup.TopMargin = 56.692913385826778
$ws2.PageSetup.BottomMargin = 56.692913385826778

# --- view state: LMR keeps focus on B28, Tabelle1 becomes the active tab -
[void]$ws1.Range("B28").Select()
[void]$ws2.Activate()
[void]$ws2.Range("E9").Select()
